# Adapt column header formatting to respective input file names.
#   *_old  -> *_FV2304
#   *_new  -> *_FV2310
# Then turn the data range into an Excel Table and freeze the header row,
# matching the target OOXML described by the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the header-row labels: "<name>_old" -> "<name>_FV2304" and
#    "<name>_new" -> "<name>_FV2310". The "diff" header in column K is
#    left untouched.
# ---------------------------------------------------------------------
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $v = $cell.Value2
    if ($v -ne $null) {
        $nv = $v -replace '_old$', '_FV2304'
        $nv = $nv -replace '_new$', '_FV2310'
        if ($nv -ne $v) {
            $cell.Value = $nv
        }
    }
}

# ---------------------------------------------------------------------
# 2. Turn A1:U64 into a proper Excel Table (ListObject) so the header
#    row gets filter buttons and a <tableParts> entry on the sheet.
#
#    The header row (row 1) already carries explicit manual formatting
#    (bold / fill / border / centered+wrap, style index "1"). If that
#    formatting is still present at the moment the table is created,
#    the engine captures it as a per-table header-style override
#    (a new <dxf> + headerRowDxfId) - which the target workbook does
#    not have. To avoid that, the existing header formatting is copied
#    aside, cleared, the table is created against the now-unformatted
#    header, and the original formatting is pasted back afterwards
#    (formatting changes made *after* table creation are not folded
#    into the table's style override).
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

$headerRange.Copy()
$ws.Range("A1000").PasteSpecial(-4122)   # xlPasteFormats

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.TableStyle = ""

$scratchRange.Copy()
$ws.Range("A1").PasteSpecial(-4122)      # xlPasteFormats
$scratchRange.Clear()

# ---------------------------------------------------------------------
# 3. Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header columns renamed, table added, header row frozen."
